$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.299.12'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').Value = '2.748.02'
$ws.Range('E3').Value = '  -5.07%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '510.88'
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.43'
$ws.Range('E6').Value = '  +1.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.539'
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('D9').Value = '2.763.12'
$ws.Range('E9').Value = '  -4.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.12'
$ws.Range('E10').Value = '  +4.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.106'
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.353'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').Value = '3.230.89'
$ws.Range('E14').Value = '  -4.60%  '
$ws.Range('D15').Value = '59.255.36'
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.05'
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.774.11'
$ws.Range('E17').Value = '  -4.16%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000138'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.81'
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.18'
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.14'
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.32'
$ws.Range('E22').Value = '  -3.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.62'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.50'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.431'
$ws.Range('E26').Value = '  -3.03%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.175'
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.994'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').Value = '0.0₃0851'
$ws.Range('E29').Value = '  +1.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.61'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.63'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.42'
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.08'
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.27'
$ws.Range('E35').Value = '  -0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.48'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.969'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('E38').Value = '  -2.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.29'
$ws.Range('E39').Value = '  -4.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.42'
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.57'
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('D42').Value = '2.202.59'
$ws.Range('E42').Value = '  -5.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0564'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.611'
$ws.Range('E44').Value = '  -4.59%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.995'
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.40'
$ws.Range('E46').Value = '  -5.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.81'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.37'
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0229'
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0893'
$ws.Range('E50').Value = '  -3.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.40'
$ws.Range('E51').Value = '  +2.14%  '
